# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Femacal de La Calera" - Mango as row 511,
# pushing the existing rows 511..542 down to 512..543.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 511 (shifts rows 511:542 down to 512:543).
$ws.Rows.Item(511).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(511, 1).Value = 3
$ws.Cells.Item(511, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(511, 3).Value = "Coquimbo"
$ws.Cells.Item(511, 4).Value = 44931
$ws.Cells.Item(511, 5).Value = 5
$ws.Cells.Item(511, 6).Value = "Fruta"
$ws.Cells.Item(511, 7).Value = 100108
$ws.Cells.Item(511, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(511, 9).Value = 100108002
$ws.Cells.Item(511, 10).Value = "Mango"
$ws.Cells.Item(511, 11).Value = "Sin especificar"
$ws.Cells.Item(511, 12).Value = "Primera"
$ws.Cells.Item(511, 13).Value = 456
$ws.Cells.Item(511, 14).Value = 7000
$ws.Cells.Item(511, 15).Value = 7500
$ws.Cells.Item(511, 16).Value = 7250
$ws.Cells.Item(511, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(511, 18).Value = "Perú"
$ws.Cells.Item(511, 19).Value = 1812
$ws.Cells.Item(511, 20).Value = 4
